# Insert a new record at row 15 (pushing the existing rows 15-49 down to 16-50),
# matching the weekly data-refresh reflected in the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15:49 down by inserting a new blank row at 15.
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Cells.Item(15, 1).Value = 2
$ws.Cells.Item(15, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(15, 3).Value = "Coquimbo"
$ws.Cells.Item(15, 4).Value = 44581
$ws.Cells.Item(15, 5).Value = 4
$ws.Cells.Item(15, 6).Value = 100112032
$ws.Cells.Item(15, 7).Value = "Zapallo italiano"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 300
$ws.Cells.Item(15, 11).Value = 11000
$ws.Cells.Item(15, 12).Value = 12000
$ws.Cells.Item(15, 13).Value = 11500
$ws.Cells.Item(15, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(15, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 16).Value = 192
$ws.Cells.Item(15, 17).Value = 60
$ws.Cells.Item(15, 18).Value = "Hortaliza"
